$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the specified rows, per the diff.
$ws.Range("F3").Value = 3
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -1
$ws.Range("F17").Value = -4
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = -4
$ws.Range("F27").Value = 3
$ws.Range("F28").Value = 0
$ws.Range("F36").Value = 2
